$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 11.450000000000001
$ws.Range("C2").Value = 4.9000000000000004
$ws.Range("D2").Value = 13.25
$ws.Range("E2").Value = 16.350000000000001

# Row 3 values
$ws.Range("B3").Value = 5.3000000000000007
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 11.15
$ws.Range("E3").ClearContents()

# Update selection to match target
$ws.Range("B1:E3").Select()
